$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The NATMI re-run added a new "ECs" sending cluster as the top result and
# shifted the previous rows down. Insert a fresh row 2 for it (push old
# rows 2-3 down to 3-4), then strip any formatting the insert copied down
# from the header row so the new row matches the plain (unstyled) data rows.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# New row 2: ECs -> Nrg2 -> Erbb4 -> MuSCs (fresh TPM numbers)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nrg2"
$ws.Range("C2").Value = "Erbb4"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1182943333333333
$ws.Range("H2").Value = 0.354883
$ws.Range("I2").Value = 0.07902913105657369
$ws.Range("J2").Value = 0.07902913105657366
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.007090666666666666
$ws.Range("N2").Value = 0.021272
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.0008387856862222222
$ws.Range("R2").Value = 0.007549071176
$ws.Range("S2").Value = 0.07902913105657369
$ws.Range("T2").Value = 0.07902913105657366

# Row 3 (was row 2: FAPs -> Nrg2 -> Erbb4 -> MuSCs) keeps its raw TPM
# numbers, but the derived-specificity columns (I, J, S, T) change because
# the specificity is computed relative to all rows, and there's a new row now.
$ws.Range("I3").Value = 0.7871529310322559
$ws.Range("J3").Value = 0.7871529310322558
$ws.Range("S3").Value = 0.7871529310322559
$ws.Range("T3").Value = 0.7871529310322558

# Row 4 (was row 3: MuSCs -> Nrg2 -> Erbb4 -> MuSCs) same story.
$ws.Range("I4").Value = 0.1338179379111705
$ws.Range("J4").Value = 0.1338179379111705
$ws.Range("S4").Value = 0.1338179379111705
$ws.Range("T4").Value = 0.1338179379111705
